# The deck has 23 slides. The two slides right before the final
# "End of Chapter" slide ("14.1 Reduction: sum, numel" title slide and its
# companion slide with the reduction-operation code walkthrough) are being
# removed entirely, leaving the "End of Chapter" slide as slide 21.
$p = $ppt.ActivePresentation

# Slide 21 = "14.1 Reduction: sum, numel" (short/title-only version)
# Slide 22 = "14.1 Reduction: sum, numel" (full version w/ code + images)
# Deleting index 21 twice removes both, since the following slide
# (the old #22) shifts down into index 21 after the first delete.
$p.Slides.Item(21).Delete()
$p.Slides.Item(21).Delete()
